$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '64.862.73'
Set-TextValue $ws.Range('D3') '3.376.67'
Set-TextValue $ws.Range('E3') '  +3.26%  '
Set-TextValue $ws.Range('E4') '  -0.14%  '
Set-TextValue $ws.Range('D5') '558.82'
Set-TextValue $ws.Range('E5') '  +4.06%  '
Set-TextValue $ws.Range('D6') '173.90'
Set-TextValue $ws.Range('E6') '  +3.12%  '
Set-TextValue $ws.Range('D7') '0.623'
Set-TextValue $ws.Range('E7') '  +2.75%  '
Set-TextValue $ws.Range('D8') '3.368.18'
Set-TextValue $ws.Range('E8') '  +3.31%  '
Set-TextValue $ws.Range('E9') '  -0.03%  '
Set-TextValue $ws.Range('D10') '0.166'
Set-TextValue $ws.Range('E10') '  +11.31%  '
Set-TextValue $ws.Range('D11') '0.631'
Set-TextValue $ws.Range('E11') '  +4.80%  '
Set-TextValue $ws.Range('D12') '53.89'
Set-TextValue $ws.Range('E12') '  +4.14%  '
Set-TextValue $ws.Range('D13') '0.0000276'
Set-TextValue $ws.Range('E13') '  +6.45%  '
Set-TextValue $ws.Range('D14') '9.09'
Set-TextValue $ws.Range('E14') '  +4.13%  '
Set-TextValue $ws.Range('D15') '3.925.65'
Set-TextValue $ws.Range('E15') '  +3.17%  '
Set-TextValue $ws.Range('D16') '18.23'
Set-TextValue $ws.Range('E16') '  +2.31%  '
Set-TextValue $ws.Range('E17') '  +3.55%  '
Set-TextValue $ws.Range('D18') '3.391.44'
Set-TextValue $ws.Range('E18') '  +3.65%  '
Set-TextValue $ws.Range('D19') '64.775.12'
Set-TextValue $ws.Range('D20') '11.78'
Set-TextValue $ws.Range('E20') '  +2.76%  '
Set-TextValue $ws.Range('D21') '0.989'
Set-TextValue $ws.Range('E21') '  +3.54%  '
Set-TextValue $ws.Range('D22') '469.72'
Set-TextValue $ws.Range('E22') '  +14.83%  '
Set-TextValue $ws.Range('E23') '  +13.53%  '
Set-TextValue $ws.Range('D24') '4.12'
Set-TextValue $ws.Range('E24') '  +4.07%  '
Set-TextValue $ws.Range('D25') '86.69'
Set-TextValue $ws.Range('E25') '  +5.62%  '
Set-TextValue $ws.Range('D26') '13.54'
Set-TextValue $ws.Range('E26') '  +2.41%  '
Set-TextValue $ws.Range('E27') '  +9.42%  '
Set-TextValue $ws.Range('D28') '10.81'
Set-TextValue $ws.Range('E28') '  +3.69%  '
Set-TextValue $ws.Range('D29') '8.73'
Set-TextValue $ws.Range('E29') '  +3.65%  '
Set-TextValue $ws.Range('D30') '30.65'
Set-TextValue $ws.Range('E30') '  +7.24%  '
Set-TextValue $ws.Range('D31') '6.75'
Set-TextValue $ws.Range('E31') '  +8.50%  '
Set-TextValue $ws.Range('D32') '11.45'
Set-TextValue $ws.Range('E32') '  +2.59%  '
Set-TextValue $ws.Range('D33') '573.13'
Set-TextValue $ws.Range('E33') '  +0.60%  '
Set-TextValue $ws.Range('D34') '61.42'
Set-TextValue $ws.Range('E34') '  +6.98%  '
Set-TextValue $ws.Range('E35') '  +3.00%  '
Set-TextValue $ws.Range('E36') '  -0.01%  '
Set-TextValue $ws.Range('D37') '3.59'
Set-TextValue $ws.Range('E37') '  +7.09%  '
Set-TextValue $ws.Range('E38') '  -2.83%  '
Set-TextValue $ws.Range('D39') '35.44'
Set-TextValue $ws.Range('E39') '  +2.50%  '
Set-TextValue $ws.Range('D40') '0.0₃0739'
Set-TextValue $ws.Range('E40') '  +2.42%  '
Set-TextValue $ws.Range('D41') '0.368'
Set-TextValue $ws.Range('E41') '  +2.80%  '
Set-TextValue $ws.Range('B42') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D42') '1.00'
Set-TextValue $ws.Range('E42') '  -0.20%  '
Set-TextValue $ws.Range('B43') 'Maker'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D43') '3.078.00'
Set-TextValue $ws.Range('E43') '  +0.18%  '
Set-TextValue $ws.Range('D44') '2.83'
Set-TextValue $ws.Range('E44') '  +5.04%  '
Set-TextValue $ws.Range('D45') '0.0414'
Set-TextValue $ws.Range('E45') '  +5.26%  '
Set-TextValue $ws.Range('D46') '0.134'
Set-TextValue $ws.Range('E46') '  +6.25%  '
Set-TextValue $ws.Range('D47') '2.45'
Set-TextValue $ws.Range('E47') '  +3.28%  '
Set-TextValue $ws.Range('D48') '3.13'
Set-TextValue $ws.Range('E48') '  -2.97%  '
Set-TextValue $ws.Range('D49') '2.61'
Set-TextValue $ws.Range('E49') '  +0.95%  '
Set-TextValue $ws.Range('D50') '137.65'
Set-TextValue $ws.Range('E50') '  +4.48%  '
Set-TextValue $ws.Range('D51') '8.28'
Set-TextValue $ws.Range('E51') '  +4.80%  '

Write-Host "Applied all changes"